$d = $word.ActiveDocument

function Find-Range([string]$needle, [int]$searchFrom) {
    $probe = $d.Range($searchFrom, $d.Content.End)
    $probe.Find.ClearFormatting()
    $found = $probe.Find.Execute($needle, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found -or -not $probe.Find.Found) {
        throw "Could not find '$needle' token"
    }
    # Re-wrap in a fresh Range: Find.Execute repositions the range's
    # Start/End, but a range that has had Find run on it confuses
    # InsertXML's insertion point (it starts appending at the story/
    # paragraph end instead of the match). A freshly minted Range with
    # the same bounds does not have that problem.
    return $d.Range($probe.Start, $probe.End)
}

function Wrap-RunsXml([string]$runsXml) {
    return '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body><w:p>' + $runsXml + '</w:p></w:body></w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'
}

# ---------------------------------------------------------------------
# 1) Split the run containing "{m" into two runs: "{" and "m".
#    InsertXML only lands correctly in place when the replaced range is
#    the LAST run of its paragraph (otherwise the new runs get appended
#    at the paragraph's end). "{m" is the first run of its paragraph, so
#    temporarily split the paragraph right after it, do the run split
#    (now valid, since "{m" is alone in its own paragraph), then merge
#    the paragraph back by deleting the inserted paragraph mark.
# ---------------------------------------------------------------------
$target = Find-Range "{m" 0
$splitStart = $target.Start
$splitEnd = $target.End
$target.InsertParagraphAfter()

$target = $d.Range($splitStart, $splitEnd)
$target.InsertXML((Wrap-RunsXml '<w:r><w:t>{</w:t></w:r><w:r><w:t>m</w:t></w:r>'))

# The paragraph mark just inserted sits right after "m" (i.e. right at
# $splitEnd, the original end of the "{m" run); delete it to merge the
# two paragraphs back into one.
$pMark = $d.Range($splitEnd, $splitEnd + 1)
$pMark.Delete()

# ---------------------------------------------------------------------
# 2) Split the run containing "')}" into two runs: "')" (keeps the
#    orange-color formatting) and "}" (plain, xml:space="preserve").
#    This run is already the last run of its paragraph, so InsertXML
#    lands in place directly.
# ---------------------------------------------------------------------
$target = Find-Range "')}" 0
$colorRPr = '<w:rPr><w:color w:themeColor="accent6" w:themeShade="BF" w:val="E36C0A"/></w:rPr>'
$quote = [string][char]39
$runsXml2 = '<w:r>' + $colorRPr + '<w:t>' + $quote + ')</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">}</w:t></w:r>'
$target.InsertXML((Wrap-RunsXml $runsXml2))
